# Generate Report for Handoff
#
# The localization status moves from "In Translation" to "Ready for
# handoff" and the associated timestamps are refreshed to the moment the
# handoff report was generated. This touches the per-locale "Status" /
# "Latest Handoff Datetime" cells on the zh-cn and de-de sheets, and the
# roll-up "Overview" sheet that mirrors the same status/date per locale.
# Widening the Status columns is a side effect of the longer new text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: E2/F2 hold the per-locale status (zh-cn / de-de),
# G2 the "Latest HO Xliff Generate Date".
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-03 15:03:39"

# zh-cn sheet: C2 status, H2 "Latest Handoff Datetime".
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-03 15:03:33"

# de-de sheet: C2 status, H2 "Latest Handoff Datetime".
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-03 15:03:39"

# The "Status" columns grew (longer new text), so they are widened to fit,
# matching the workbook's recorded autofit/resize of those columns.
$newStatusWidth = 16.333333333333332
$wsOverview.Columns.Item(5).ColumnWidth = $newStatusWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $newStatusWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newStatusWidth
